# day3/미리가는 연구실(07.24).pptx - slide 3, "TextBox 2"
# Paragraph 2 reads: "가상환경 실행 conda activate Arduino"
# The commit lower-cases the trailing "Arduino" -> "arduino"
# (PowerPoint records this as a run split: " activate " / "arduino").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item("TextBox 2")
$tr = $shp.TextFrame.TextRange

# Locate paragraph 2 ("가상환경 실행 conda activate Arduino") and, within it,
# the "Arduino" substring, then retype it in lower case.
$para = $tr.Paragraphs(2, 1)
$relIdx = $para.Text.IndexOf("Arduino")
$absStart = $para.Start + $relIdx

$target = $tr.Characters($absStart, 7)
$target.Text = "arduino"
